$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the student name looked up in the VLOOKUP example (G16):
# change from "Manuel" to "Raúl"
$ws.Range("G16").Value = "Raúl"

# Move the active selection to G21, matching where the user clicked next
$ws.Range("G21").Select()
